$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4029636681079865
$ws.Range("B1").Value = 0.660445511341095
$ws.Range("C1").Value = 2.069004058837891
$ws.Range("D1").Value = 4.81011438369751
$ws.Range("E1").Value = 2.126118183135986
